# Borrar puntos al final de frases en viñetas para versionas Full CV
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# Quitar el punto final de la frase del Ig Nobel Prize (columna E, fila 2)
$ws.Range("E2").Value = "Por ‘tratar de cuantificar la relación entre la desigualdad de ingresos nacionales en diferentes países y la cantidad promedio de besos boca a boca’ (Watkins,, et al., 2019)"

# Quitar el punto final de la frase del Annual Prize in Evolutionary Psychology (columna E, fila 8)
$ws.Range("E8").Value = "Mejor desempeño general en la maestría"

# Mover la selección activa a E12, como quedó tras la edición
$ws.Range("E12").Select()
